$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Preston")

# --- Add Preston's hours from the initial project discussion meetings ---

# Row 2 already has a date-formatted style (s=5) on A2; copy that format
# down to A3:A5 before filling in the values so the new rows match the
# existing "date" column formatting exactly.
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)

# Columns B (hours) and C (task) inherit the plain column style, but copy
# explicitly from the existing row 2 cells to stay consistent.
$ws.Range("B2").Copy()
$ws.Range("B3:B5").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C3:C5").PasteSpecial(-4122)

$ws.Range("A2").Value = 45531
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Initial team meeting, discussed project and advisor ideas"

$ws.Range("A3").Value = 45532
$ws.Range("B3").Value = 0.75
$ws.Range("C3").Value = "Met with Dr. Gallagher to discuss Neromorphic Constraint Satisfaction project, and discussed with team"

$ws.Range("A4").Value = 45533
$ws.Range("B4").Value = 0.75
$ws.Range("C4").Value = "Met with Jeremy Hill to discuss TA scheduling software (generally, scheduling with constraints)"

$ws.Range("A5").Value = 45533
$ws.Range("B5").Value = 0.75
$ws.Range("C5").Value = "Met with Dr. Abuaitah to discuss assembly code simulator/educational tool project, and discussed with team"

# Widen the task column so the longer notes are readable.
$ws.Columns.Item(3).ColumnWidth = 85.15

# Preston's sheet becomes the active tab/selection (previously it was Izzy's).
$ws.Activate()
[void]$ws.Range("B6").Select()

Write-Output "done"
